# Queue-animation deck: tweak the "not enough time before job 4" caption
# that appears (identically) on the three consecutive build/animation
# slides, shortening it to "insufficient time before 4".
#
# On the first two occurrences (slides 6 and 7) the whole run's text is
# simply replaced. On the third occurrence (slide 8) the same caption is
# retyped as two runs: "insufficient time before " + "4", mirroring how
# the source deck was actually edited there.

$p = $ppt.ActivePresentation

$oldText = "not enough time before job 4"
$newText = "insufficient time before 4"

# --- Slides 6 and 7: straightforward single-run text replacement -----
# Use TextRange.Replace (rather than assigning .Text on the whole
# paragraph) so the existing run is updated in place instead of being
# split apart by a prefix/suffix text-diff.
foreach ($slideIndex in 6, 7) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(1)
    $textRange = $shape.TextFrame.TextRange
    [void]$textRange.Replace($oldText, $newText)
}

# --- Slide 8: same wording change, but split into two runs ("...before " / "4") --
$slide8 = $p.Slides.Item(8)
$shape8 = $slide8.Shapes.Item(1)
$textRange8 = $shape8.TextFrame.TextRange
$paragraph8 = $textRange8.Paragraphs(7, 1)

if ($paragraph8.Text.TrimEnd("`r") -eq $oldText) {
    # Keep the trailing "4" as its own (untouched) run, and retype
    # everything before it as a new run reading "insufficient time before ".
    $prefix = $paragraph8.Characters(1, $paragraph8.Length - 2)
    $prefix.Text = "insufficient time before "
}
